$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append "(N" (instruction ordinal, 1-based) to each instruction-name cell
# in A34:A61 -- e.g. "Add" -> "Add(1", "Bltz" -> "Bltz(28".
for ($row = 34; $row -le 61; $row++) {
    $n = $row - 33
    $cell = $ws.Cells.Item($row, 1)
    $current = $cell.Value()
    $cell.Value = "$current($n"
}

# Update the active selection shown in the saved sheet view.
$ws.Range("N54").Select() | Out-Null
